$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2..41 of columns A (KILLS) and E (DEATHS) were stored as inline-string
# text ("0", "1", "2", ...) even though they are numeric counters. Re-enter
# them as real numbers so the cells become numeric (t="n") instead of text.
# Two values were also genuinely wrong and get corrected at the same time:
#   - A20 held the stray letter "S" and should read 5 (kills keep climbing by
#     frame, same as the surrounding rows).
#   - A41/E41/F41 all held the placeholder text "erro"; the real kill/death/
#     assist counts for that frame are 10, 1 and 7 respectively.
$kills  = @(0,0,0,0,0,1,1,1,1,2,2,2,3,3,3,3,3,5,5,5,5,6,6,7,7,7,7,7,7,8,8,8,8,8,8,9,9,10,10,10)
$deaths = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)

for ($i = 0; $i -lt $kills.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $kills[$i]
    $ws.Range("E$row").Value = $deaths[$i]
}

# Column F keeps its text type everywhere; only row 41's placeholder changes
# value ("erro" -> "7"). The leading apostrophe forces Excel to store the
# digit as text instead of re-interpreting it as a number, and resetting the
# style afterwards avoids leaving a stray "quote prefix" cell format behind.
$ws.Range("F41").Value = "'7"
$ws.Range("F41").Style = "Normal"
